$wb = $excel.ActiveWorkbook

# --- "usuario" sheet (sheet7.xml): insert two new attribute rows -----------
$wsUsuario = $wb.Worksheets.Item("usuario")

# Insert two blank rows above the current row 3 ("id_rol" FK row), shifting
# the existing rows 3-4 down to 5-6.
$wsUsuario.Rows.Item(3).Insert()
$wsUsuario.Rows.Item(3).Insert()

# New row 3: nombre / VARCHAR(50) / NOT NULL / (no LLAVE) / nombre del usuario
$wsUsuario.Range("A3").Value = "nombre"
$wsUsuario.Range("B3").Value = "VARCHAR(50)"
$wsUsuario.Range("C3").Value = "NOT NULL"
$wsUsuario.Range("E3").Value = "nombre del usuario"

# New row 4: apellido / VARCHAR(50) / NOT NULL / (no LLAVE) / apellido del usuario
$wsUsuario.Range("A4").Value = "apellido"
$wsUsuario.Range("B4").Value = "VARCHAR(50)"
$wsUsuario.Range("C4").Value = "NOT NULL"
$wsUsuario.Range("E4").Value = "apellido del usuario"

# --- sheet selection / active-cell bookkeeping -----------------------------
# "estado_aut" (sheet8.xml) loses the tab selection and its remembered
# selection moves to E9.
$wsEstadoAut = $wb.Worksheets.Item("estado_aut")
$null = $wsEstadoAut.Range("E9").Select()

# "usuario" (sheet7.xml) becomes the active tab, with E5 selected (the cell
# that used to be E3 before the two rows were inserted).
$null = $wsUsuario.Range("E5").Select()
